# Apply the "output generated at 456a3b4" gh-pages data refresh.
# Workbook has 4 sheets: 1=展览(Exhibition) 2=演出(Performance)
# 3=本地生活(Local Life) 4=全部类型(All Types).
# Most edits are simple numeric refreshes of column F ("views"/interest
# counter). Sheet 4 additionally has rows 11-18 replaced: the old
# 2024-07-06 concert row is removed, rows 12-18 shift up to 11-17, and a
# brand-new row (2024-07-19 萤火虫动漫游戏嘉年华...) lands in the vacated
# slot 18.

$wb = $excel.ActiveWorkbook

function Set-TextCell {
    param($ws, [string]$ref, [string]$text)
    $ws.Range($ref).NumberFormat = "@"
    $ws.Range($ref).Value = $text
}

function Set-NumCell {
    param($ws, [string]$ref, $num)
    $ws.Range($ref).Value = $num
}

# ---------------------------------------------------------------
# Sheet 1: 展览 (Exhibition) -- column F refreshes only
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
Set-NumCell $ws1 "F3"  1236
Set-NumCell $ws1 "F4"  1299
Set-NumCell $ws1 "F6"  182
Set-NumCell $ws1 "F7"  564
Set-NumCell $ws1 "F8"  32
Set-NumCell $ws1 "F9"  359
Set-NumCell $ws1 "F11" 1287
Set-NumCell $ws1 "F12" 29474
Set-NumCell $ws1 "F13" 4968
Set-NumCell $ws1 "F14" 48
Set-NumCell $ws1 "F15" 269
Set-NumCell $ws1 "F16" 496
Set-NumCell $ws1 "F19" 30
Set-NumCell $ws1 "F21" 354
Set-NumCell $ws1 "F22" 19
Set-NumCell $ws1 "F24" 284
Set-NumCell $ws1 "F26" 363
Set-NumCell $ws1 "F28" 90
Set-NumCell $ws1 "F30" 676
Set-NumCell $ws1 "F31" 222
Set-NumCell $ws1 "F33" 569
Set-NumCell $ws1 "F34" 81
Set-NumCell $ws1 "F36" 657
Set-NumCell $ws1 "F39" 10

# ---------------------------------------------------------------
# Sheet 2: 演出 (Performance) -- column F refreshes only
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)
Set-NumCell $ws2 "F7"  926
Set-NumCell $ws2 "F12" 4272
Set-NumCell $ws2 "F23" 4257

# ---------------------------------------------------------------
# Sheet 3: 本地生活 (Local Life) -- column F refresh only
# ---------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)
Set-NumCell $ws3 "F4" 1237

# ---------------------------------------------------------------
# Sheet 4: 全部类型 (All Types)
# ---------------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)

# Simple column F refreshes (rows untouched by the shift below)
Set-NumCell $ws4 "F4"  1237
Set-NumCell $ws4 "F9"  926
Set-NumCell $ws4 "F10" 1236
Set-NumCell $ws4 "F26" 269
Set-NumCell $ws4 "F29" 496
Set-NumCell $ws4 "F31" 30
Set-NumCell $ws4 "F35" 19
Set-NumCell $ws4 "F37" 284
Set-NumCell $ws4 "F39" 90
Set-NumCell $ws4 "F41" 676
Set-NumCell $ws4 "F43" 222
Set-NumCell $ws4 "F47" 81
Set-NumCell $ws4 "F49" 657

# Row 11 ("2024-07-06 龙珠/灌篮高手 音乐会") is removed from the listing.
# Rows 12-18 shift up to become rows 11-17 (their own F values also tick
# up as part of the same refresh), and a freshly-scraped row takes the
# vacated slot 18.

# New row 11 (was row 12): 特摄FansMeetup
Set-TextCell $ws4 "B11" "2024-07-13"
Set-TextCell $ws4 "C11" "广州·特摄FansMeetup"
Set-TextCell $ws4 "D11" "芳村大道下市直街1号信义会馆21栋(近白鹅潭风情酒吧街) 信义会馆-21栋"
Set-TextCell $ws4 "E11" "2024.07.13 10:00-07.13 19:00"
Set-NumCell  $ws4 "F11" 182
Set-NumCell  $ws4 "G11" 69.90000000000001
Set-TextCell $ws4 "H11" "https://show.bilibili.com/platform/detail.html?id=87031"
Set-TextCell $ws4 "I11" "//i1.hdslb.com/bfs/openplatform/202406/9ffC9a8n1717578946827.jpeg"

# New row 12 (was row 13): 第5人格only3.0联动特别篇
Set-TextCell $ws4 "B12" "2024-07-13"
Set-TextCell $ws4 "C12" "广州·第5人格only3.0联动特别篇"
Set-TextCell $ws4 "D12" "奥体南路12号 优托邦(奥体旗舰店)"
Set-TextCell $ws4 "E12" "2024.07.13 10:00-07.13 17:00"
Set-NumCell  $ws4 "F12" 564
Set-NumCell  $ws4 "G12" 54
Set-TextCell $ws4 "H12" "https://show.bilibili.com/platform/detail.html?id=86740"
Set-TextCell $ws4 "I12" "//i2.hdslb.com/bfs/openplatform/202405/mwlJqj0o1717149700846.jpeg"

# New row 13 (was row 14): AI动漫嘉年华6.0-原神coser免票
Set-TextCell $ws4 "B13" "2024-07-14"
Set-TextCell $ws4 "C13" "广州·AI动漫嘉年华6.0-原神coser免票"
Set-TextCell $ws4 "D13" "奥体南路12号 优托邦(奥体旗舰店)"
Set-TextCell $ws4 "E13" "2024.07.14 10:00-07.14 17:00"
Set-NumCell  $ws4 "F13" 32
Set-NumCell  $ws4 "G13" 45
Set-TextCell $ws4 "H13" "https://show.bilibili.com/platform/detail.html?id=88120"
Set-TextCell $ws4 "I13" "//i2.hdslb.com/bfs/openplatform/202406/0GnShKZT1719367360057.jpeg"

# New row 14 (was row 15): OOPS 7th
Set-TextCell $ws4 "B14" "2024-07-14"
Set-TextCell $ws4 "C14" "广州·OOPS 7th"
Set-TextCell $ws4 "D14" "流花街道流花路119号越秀公园站B2、C出口 广州越秀国际会议中心"
Set-TextCell $ws4 "E14" "2024.07.14 09:30-07.15 17:00"
Set-NumCell  $ws4 "F14" 359
Set-NumCell  $ws4 "G14" 20
Set-TextCell $ws4 "H14" "https://show.bilibili.com/platform/detail.html?id=87550"
Set-TextCell $ws4 "I14" "//i2.hdslb.com/bfs/openplatform/202405/Qi8gB0Bi1715922859908.png"

# New row 15 (was row 16): 《大鱼.刀剑如梦》最美国风经典影视金曲音乐会
Set-TextCell $ws4 "B15" "2024-07-14"
Set-TextCell $ws4 "C15" "广州·《大鱼.刀剑如梦》最美国风经典影视金曲音乐会"
Set-TextCell $ws4 "D15" "东风中路299号 广州中山纪念堂"
Set-TextCell $ws4 "E15" "2024.07.14 19:30-07.14 21:30"
Set-NumCell  $ws4 "F15" 0
Set-NumCell  $ws4 "G15" 85
Set-TextCell $ws4 "H15" "https://show.bilibili.com/platform/detail.html?id=87899"
Set-TextCell $ws4 "I15" "//i0.hdslb.com/bfs/openplatform/202406/WJUHqwHD1718878927800.png"

# New row 16 (was row 17): 幻毛纪AnimalFurryOnly
Set-TextCell $ws4 "B16" "2024-07-14"
Set-TextCell $ws4 "C16" "广州·幻毛纪AnimalFurryOnly"
Set-TextCell $ws4 "D16" "芳村大道下市直街1号信义会馆21栋(近白鹅潭风情酒吧街) 信义会馆-21栋"
Set-TextCell $ws4 "E16" "2024.07.14 10:00-07.14 19:00"
Set-NumCell  $ws4 "F16" 55
Set-NumCell  $ws4 "G16" 68.8
Set-TextCell $ws4 "H16" "https://show.bilibili.com/platform/detail.html?id=87273"
Set-TextCell $ws4 "I16" "//i0.hdslb.com/bfs/openplatform/202406/9z1DMHsl1718181280279.png"

# New row 17 (was row 18): 火影only
Set-TextCell $ws4 "B17" "2024-07-14"
Set-TextCell $ws4 "C17" "广州·火影only"
Set-TextCell $ws4 "D17" "人和镇蚌湖清河大街168号 人和园"
Set-TextCell $ws4 "E17" "2024.07.14 09:30-07.14 17:30"
Set-NumCell  $ws4 "F17" 1287
Set-NumCell  $ws4 "G17" 78
Set-TextCell $ws4 "H17" "https://show.bilibili.com/platform/detail.html?id=84815"
Set-TextCell $ws4 "I17" "//i2.hdslb.com/bfs/openplatform/202404/QLOhW4Nr1714384036670.png"

# New row 18 (brand new): 萤火虫动漫游戏嘉年华 x KKWORLD2024 快看漫画乐园
# Note column G switches from a numeric price to the text "已售罄" (sold out).
Set-TextCell $ws4 "B18" "2024-07-19"
Set-TextCell $ws4 "C18" "广州·萤火虫动漫游戏嘉年华 × KKWORLD2024 快看漫画乐园"
Set-TextCell $ws4 "D18" "新港东路1000号 保利世贸博览馆"
Set-TextCell $ws4 "E18" "2024.07.19 09:00-07.22 17:00"
Set-NumCell  $ws4 "F18" 29474
Set-TextCell $ws4 "G18" "已售罄"
Set-TextCell $ws4 "H18" "https://show.bilibili.com/platform/detail.html?id=87210"
Set-TextCell $ws4 "I18" "//i1.hdslb.com/bfs/openplatform/202406/DTCdOTPs1718177177472.jpeg"
